$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the data set. It belongs logically right
# after the existing row 310, so insert a fresh row at 311 (this pushes the
# former rows 311-377 down to 312-378, which is exactly what the target
# workbook shows: dimension grows from A1:T377 to A1:T378).
$ws.Rows.Item(311).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Cells.Item(311, 1).Value  = 9
$ws.Cells.Item(311, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(311, 3).Value  = "Metropolitana"
$ws.Cells.Item(311, 4).Value  = 44641
$ws.Cells.Item(311, 5).Value  = 13
$ws.Cells.Item(311, 6).Value  = "Fruta"
$ws.Cells.Item(311, 7).Value  = 100108
$ws.Cells.Item(311, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(311, 9).Value  = 100108002
$ws.Cells.Item(311, 10).Value = "Mango"
$ws.Cells.Item(311, 11).Value = "Sin especificar"
$ws.Cells.Item(311, 12).Value = "Primera"
$ws.Cells.Item(311, 13).Value = 580
$ws.Cells.Item(311, 14).Value = 7500
$ws.Cells.Item(311, 15).Value = 8000
$ws.Cells.Item(311, 16).Value = 7741
$ws.Cells.Item(311, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(311, 18).Value = "Perú"
$ws.Cells.Item(311, 19).Value = 1935
$ws.Cells.Item(311, 20).Value = 4

# Keep the date column's style consistent with the rest of column D.
$ws.Cells.Item(311, 4).NumberFormat = $ws.Cells.Item(312, 4).NumberFormat
